# [Fonds de solidarite] Add 2021-01-22 data
#
# Updates nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# figures for a handful of region / classe_effectif rows to reflect the
# 2021-01-22 data refresh. Source values are text (t="inlineStr") in the
# original workbook, so we re-enter them with a leading apostrophe to force
# Excel to keep them as text instead of auto-converting to numbers, then
# reset the cell style back to "Normal" so no stray number-format / style
# gets attached to the cell (matches the original formatting exactly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row 4 - Auvergne-Rhône-Alpes / 3 à 5 salariés
Set-TextValue $ws "C4" "671"
Set-TextValue $ws "D4" "566"
Set-TextValue $ws "E4" "6908275.49"

# Row 25 - Centre-Val de Loire / 6 à 9 salariés
Set-TextValue $ws "C25" "85"
Set-TextValue $ws "E25" "1379797.91"

# Row 37 - Grand Est / 3 à 5 salariés
Set-TextValue $ws "C37" "386"
Set-TextValue $ws "E37" "3350933.69"

# Row 62 - Île-de-France / 0 salarié
Set-TextValue $ws "C62" "1140"
Set-TextValue $ws "E62" "3803659.28"

# Row 63 - Île-de-France / 1 ou 2 salariés
Set-TextValue $ws "C63" "5716"
Set-TextValue $ws "E63" "24531453.78"

# Row 64 - Île-de-France / 3 à 5 salariés
Set-TextValue $ws "C64" "3149"
Set-TextValue $ws "E64" "19489380.84"

# Row 65 - Île-de-France / 6 à 9 salariés
Set-TextValue $ws "C65" "1116"
Set-TextValue $ws "D65" "1012"
Set-TextValue $ws "E65" "9135015.17"

# Row 66 - Île-de-France / 10 à 19 salariés
Set-TextValue $ws "C66" "317"
Set-TextValue $ws "D66" "286"
Set-TextValue $ws "E66" "4380004.41"

# Row 67 - Île-de-France / 20 à 49 salariés
Set-TextValue $ws "C67" "51"
Set-TextValue $ws "E67" "1887881.18"

# Row 94 - Nouvelle-Aquitaine / 3 à 5 salariés
Set-TextValue $ws "C94" "488"
Set-TextValue $ws "E94" "4743869.46"

Write-Output "Updated 10 rows (4, 25, 37, 62-67, 94) with 2021-01-22 data"
